$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet named "Sheet2" positioned after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with data
$ws2.Range("A1").Value = "Company"
$ws2.Range("B1").Value = "Location"
$ws2.Range("A2").Value = "Tata"
$ws2.Range("B2").Value = "Ambarnath"
$ws2.Range("A3").Value = "Mahindra"
$ws2.Range("B3").Value = "Thane"
$ws2.Range("A4").Value = "Shine"
$ws2.Range("B4").Value = "Badlapur"

# Make Sheet2 the active sheet / selected tab
$ws2.Select()
$ws2.Range("B4").Select()
